$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3549.1667
$ws.Range("I40").Value = 3380
$ws.Range("J40").Value = 3593.6843
$ws.Range("K40").Value = 3380
$ws.Range("L40").Value = 3593.6843
$ws.Range("M40").Value = -3205
$ws.Range("N40").Value = -3943.6843
$ws.Range("H64").Value = 3684.7964
$ws.Range("I64").Value = 3488.889
$ws.Range("J64").Value = 3880.7036
$ws.Range("K64").Value = 3488.889
$ws.Range("L64").Value = 3880.7036
$ws.Range("M64").Value = -3240.889
$ws.Range("N64").Value = -4376.7036
$ws.Range("H67").Value = 3684.7964
$ws.Range("I67").Value = 3488.889
$ws.Range("J67").Value = 3880.7036
$ws.Range("K67").Value = 3488.889
$ws.Range("L67").Value = 3880.7036
$ws.Range("M67").Value = -2630.889
$ws.Range("N67").Value = -5596.7036
$ws.Range("H76").Value = 3061
$ws.Range("I76").Value = 2622
$ws.Range("J76").Value = 3500
$ws.Range("K76").Value = 2622
$ws.Range("L76").Value = 3500
$ws.Range("M76").Value = -2307
$ws.Range("N76").Value = -4130
$ws.Range("H79").Value = 3061
$ws.Range("I79").Value = 2622
$ws.Range("J79").Value = 3500
$ws.Range("K79").Value = 2622
$ws.Range("L79").Value = 3500
$ws.Range("M79").Value = -1530
$ws.Range("N79").Value = -5684
$ws.Range("H136").Value = 33647.5
$ws.Range("J136").Value = 33647.5
$ws.Range("L136").Value = 33647.5
$ws.Range("N136").Value = -43847.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1520.8077
$ws.Range("I2").Value = 1472
$ws.Range("J2").Value = 2106.5
$ws.Range("K2").Value = 1472
$ws.Range("L2").Value = 2106.5
$ws.Range("M2").Value = -1359
$ws.Range("N2").Value = -2332.5
$ws.Range("H32").Value = 15397353
$ws.Range("I32").Value = 20412510
$ws.Range("K32").Value = 20412510
$ws.Range("M32").Value = -20412223
$ws.Range("H45").Value = 1111.64
$ws.Range("I45").Value = 809.95
$ws.Range("J45").Value = 2318.4
$ws.Range("K45").Value = 809.95
$ws.Range("L45").Value = 2318.4
$ws.Range("M45").Value = -432.95
$ws.Range("N45").Value = -3072.4
$ws.Range("H116").Value = 1520.8077
$ws.Range("I116").Value = 1472
$ws.Range("J116").Value = 2106.5
$ws.Range("K116").Value = 1472
$ws.Range("L116").Value = 2106.5
$ws.Range("M116").Value = 822
$ws.Range("N116").Value = -6694.5
$ws.Range("H139").Value = 39800
$ws.Range("J139").Value = 39800
$ws.Range("L139").Value = 39800
$ws.Range("N139").Value = -50080
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1520.8077
$ws.Range("I3").Value = 1472
$ws.Range("J3").Value = 2106.5
$ws.Range("K3").Value = 1472
$ws.Range("L3").Value = 2106.5
$ws.Range("M3").Value = -1358
$ws.Range("N3").Value = -2334.5
$ws.Range("H107").Value = 1887.9565
$ws.Range("I107").Value = 1874.5349
$ws.Range("J107").Value = 2080.3333
$ws.Range("K107").Value = 1874.5349
$ws.Range("L107").Value = 2080.3333
$ws.Range("M107").Value = 45.46509999999989
$ws.Range("N107").Value = -5920.3333
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2475.5
$ws.Range("I16").Value = 1248.3334
$ws.Range("J16").Value = 3211.8
$ws.Range("K16").Value = 1248.3334
$ws.Range("L16").Value = 3211.8
$ws.Range("M16").Value = -961.3334
$ws.Range("N16").Value = -3785.8
$ws.Range("H31").Value = 3319.8064
$ws.Range("I31").Value = 2240.2727
$ws.Range("J31").Value = 5958.6665
$ws.Range("K31").Value = 2240.2727
$ws.Range("L31").Value = 5958.6665
$ws.Range("M31").Value = -1945.2727
$ws.Range("N31").Value = -6548.6665
$ws.Range("H34").Value = 3319.8064
$ws.Range("I34").Value = 2240.2727
$ws.Range("J34").Value = 5958.6665
$ws.Range("K34").Value = 2240.2727
$ws.Range("L34").Value = 5958.6665
$ws.Range("M34").Value = -2038.2727
$ws.Range("N34").Value = -6362.6665
$ws.Range("H44").Value = 6875
$ws.Range("J44").Value = 6875
$ws.Range("L44").Value = 6875
$ws.Range("N44").Value = -7759
$ws.Range("H62").Value = 3161.5715
$ws.Range("I62").Value = 3000
$ws.Range("J62").Value = 3323.1428
$ws.Range("K62").Value = 3000
$ws.Range("L62").Value = 3323.1428
$ws.Range("M62").Value = -2376
$ws.Range("N62").Value = -4571.1428
$ws.Range("H65").Value = 3161.5715
$ws.Range("I65").Value = 3000
$ws.Range("J65").Value = 3323.1428
$ws.Range("K65").Value = 15000
$ws.Range("L65").Value = 16615.714
$ws.Range("M65").Value = -11880
$ws.Range("N65").Value = -22855.714
$ws.Range("H113").Value = 2475.5
$ws.Range("I113").Value = 1248.3334
$ws.Range("J113").Value = 3211.8
$ws.Range("K113").Value = 1248.3334
$ws.Range("L113").Value = 3211.8
$ws.Range("M113").Value = 921.6666
$ws.Range("N113").Value = -7551.8
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 719.2353000000001
$ws.Range("I131").Value = 503.85715
$ws.Range("K131").Value = 1511.57145
$ws.Range("M131").Value = 3528.42855
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 17555.416
$ws.Range("I113").Value = 2000
$ws.Range("J113").Value = 18969.545
$ws.Range("K113").Value = 2000
$ws.Range("L113").Value = 18969.545
$ws.Range("M113").Value = 170
$ws.Range("N113").Value = -23309.545
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4777.4
$ws.Range("I40").Value = 3819.077
$ws.Range("J40").Value = 6557.143
$ws.Range("K40").Value = 3819.077
$ws.Range("L40").Value = 6557.143
$ws.Range("M40").Value = -3683.077
$ws.Range("N40").Value = -6829.143
$ws.Range("H122").Value = 3996.9348
$ws.Range("I122").Value = 4052.6428
$ws.Range("J122").Value = 3910.2778
$ws.Range("K122").Value = 12157.9284
$ws.Range("L122").Value = 11730.8334
$ws.Range("M122").Value = -9707.928400000001
$ws.Range("N122").Value = -16630.8334
$ws.Range("H140").Value = 54214.5
$ws.Range("J140").Value = 54214.5
$ws.Range("L140").Value = 54214.5
$ws.Range("N140").Value = -64574.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H140").Value = 29214.5
$ws.Range("J140").Value = 29214.5
$ws.Range("L140").Value = 29214.5
$ws.Range("N140").Value = -39574.5
